$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.44416694687935
$ws.Range("C2").Value = 12.34343551212329
$ws.Range("D2").Value = 5.983537757823745
$ws.Range("E2").Value = 12.97702397613127
$ws.Range("G2").Value = 29.31931569586905
$ws.Range("H2").Value = 14.40973664803156
$ws.Range("L2").Value = 9.858098002935328
$ws.Range("M2").Value = 14.62578076353742
$ws.Range("O2").Value = 22.01437238511651
$ws.Range("B3").Value = 14.86894698773107
$ws.Range("C3").Value = 12.16682646210311
$ws.Range("D3").Value = 5.864507215201311
$ws.Range("E3").Value = 13.03395470205515
$ws.Range("G3").Value = 29.39768940359202
$ws.Range("H3").Value = 14.47020681464048
$ws.Range("L3").Value = 9.865552672570528
$ws.Range("M3").Value = 14.4976060266957
$ws.Range("O3").Value = 22.10867947217391
$ws.Range("B4").Value = 14.50519841127793
$ws.Range("C4").Value = 12.05707589773945
$ws.Range("D4").Value = 5.791969271045402
$ws.Range("E4").Value = 13.07073553098537
$ws.Range("G4").Value = 29.45847696110453
$ws.Range("H4").Value = 14.51036322126207
$ws.Range("L4").Value = 9.871474638594021
$ws.Range("M4").Value = 14.4200589478296
$ws.Range("O4").Value = 22.17287041197277
$ws.Range("B5").Value = 14.35452226847344
$ws.Range("C5").Value = 12.0120604122423
$ws.Range("D5").Value = 5.76259292998526
$ws.Range("E5").Value = 13.08618431730515
$ws.Range("G5").Value = 29.48641307433773
$ws.Range("H5").Value = 14.52748749962911
$ws.Range("L5").Value = 9.874226508035587
$ws.Range("M5").Value = 14.38877428396388
$ws.Range("O5").Value = 22.2006032095446
$ws.Range("B6").Value = 14.32936147948111
$ws.Range("C6").Value = 12.00456917856596
$ws.Range("D6").Value = 5.757727439260575
$ws.Range("E6").Value = 13.08877741736394
$ws.Range("G6").Value = 29.49124235458087
$ws.Range("H6").Value = 14.53037685706614
$ws.Range("L6").Value = 9.874703918052171
$ws.Range("M6").Value = 14.3835993638325
$ws.Range("O6").Value = 22.205303129393
$ws.Range("B7").Value = 14.50317594836396
$ws.Range("C7").Value = 12.05646993078742
$ws.Range("D7").Value = 5.791572287715761
$ws.Range("E7").Value = 13.07094201318301
$ws.Range("G7").Value = 29.45884093181582
$ws.Range("H7").Value = 14.51059108860513
$ws.Range("L7").Value = 9.871510379638035
$ws.Range("M7").Value = 14.41963571617739
$ws.Range("O7").Value = 22.17323805954518
$ws.Range("B8").Value = 15.24814220634168
$ws.Range("C8").Value = 12.28283593775905
$ws.Range("D8").Value = 5.942409773138172
$ws.Range("E8").Value = 12.99627570201718
$ws.Range("G8").Value = 29.34369861132362
$ws.Range("H8").Value = 14.42995776780122
$ws.Range("L8").Value = 9.860389564230598
$ws.Range("M8").Value = 14.58136294076206
$ws.Range("O8").Value = 22.04558112620694
$ws.Range("B9").Value = 16.61710646585573
$ws.Range("C9").Value = 12.7146763887369
$ws.Range("D9").Value = 6.240564768325334
$ws.Range("E9").Value = 12.86427475654835
$ws.Range("G9").Value = 29.2192080891677
$ws.Range("H9").Value = 14.29590532233552
$ws.Range("L9").Value = 9.849229591916489
$ws.Range("M9").Value = 14.90646941252111
$ws.Range("O9").Value = 21.84539358293511
$ws.Range("B10").Value = 17.55745042859276
$ws.Range("C10").Value = 13.0224554734494
$ws.Range("D10").Value = 6.458488985148474
$ws.Range("E10").Value = 12.77599673575592
$ws.Range("G10").Value = 29.19045288884514
$ws.Range("H10").Value = 14.21215506991638
$ws.Range("L10").Value = 9.847489198705009
$ws.Range("M10").Value = 15.14856791521595
$ws.Range("O10").Value = 21.72924702062728
$ws.Range("B11").Value = 17.96942524486289
$ws.Range("C11").Value = 13.15998071546784
$ws.Range("D11").Value = 6.556870121704637
$ws.Range("E11").Value = 12.73770795912675
$ws.Range("G11").Value = 29.19114094962651
$ws.Range("H11").Value = 14.17726914536837
$ws.Range("L11").Value = 9.8480916213594
$ws.Range("M11").Value = 15.25906350536267
$ws.Range("O11").Value = 21.68319870020852
$ws.Range("B12").Value = 18.123047475284
$ws.Range("C12").Value = 13.2116664186296
$ws.Range("D12").Value = 6.59397609361738
$ws.Range("E12").Value = 12.72347640599931
$ws.Range("G12").Value = 29.19338996910205
$ws.Range("H12").Value = 14.16452191485488
$ws.Range("L12").Value = 9.848519368167041
$ws.Range("M12").Value = 15.30092861220415
$ws.Range("O12").Value = 21.66674313293061
$ws.Range("B13").Value = 18.09006982984396
$ws.Range("C13").Value = 13.2005529550027
$ws.Range("D13").Value = 6.585991897914345
$ws.Range("E13").Value = 12.72652954462589
$ws.Range("G13").Value = 29.19281705840599
$ws.Range("H13").Value = 14.16724663146215
$ws.Range("L13").Value = 9.848418381389498
$ws.Range("M13").Value = 15.29191166416089
$ws.Range("O13").Value = 21.67024337037701
$ws.Range("B14").Value = 17.98211220308442
$ws.Range("C14").Value = 13.16424098189027
$ws.Range("D14").Value = 6.559926058781565
$ws.Range("E14").Value = 12.73653176509693
$ws.Range("G14").Value = 29.19128608691589
$ws.Range("H14").Value = 14.17621112979785
$ws.Range("L14").Value = 9.848122816778407
$ws.Range("M14").Value = 15.26250748350119
$ws.Range("O14").Value = 21.68182517834919
$ws.Range("B15").Value = 17.91567152496223
$ws.Range("C15").Value = 13.14194677772688
$ws.Range("D15").Value = 6.543939435608642
$ws.Range("E15").Value = 12.74269322377291
$ws.Range("G15").Value = 29.19060748582788
$ws.Range("H15").Value = 14.18176252109988
$ws.Range("L15").Value = 9.847967745009049
$ws.Range("M15").Value = 15.24449867516001
$ws.Range("O15").Value = 21.68904741831307
$ws.Range("B16").Value = 17.5301991317848
$ws.Range("C16").Value = 13.01341486330464
$ws.Range("D16").Value = 6.4520408716854
$ws.Range("E16").Value = 12.77853650758221
$ws.Range("G16").Value = 29.1906857778522
$ws.Range("H16").Value = 14.21449970690348
$ws.Range("L16").Value = 9.847477811452979
$ws.Range("M16").Value = 15.1413515632237
$ws.Range("O16").Value = 21.73239347689925
$ws.Range("B17").Value = 17.28959433162447
$ws.Range("C17").Value = 12.93390300630691
$ws.Range("D17").Value = 6.395441559248763
$ws.Range("E17").Value = 12.80100310296595
$ws.Range("G17").Value = 29.19426727359529
$ws.Range("H17").Value = 14.23540667336973
$ws.Range("L17").Value = 9.847533742990713
$ws.Range("M17").Value = 15.07814627979445
$ws.Range("O17").Value = 21.7607275682999
$ws.Range("B18").Value = 17.14972407015206
$ws.Range("C18").Value = 12.88793797592103
$ws.Range("D18").Value = 6.362818474632141
$ws.Range("E18").Value = 12.81410131837218
$ws.Range("G18").Value = 29.19762291350317
$ws.Range("H18").Value = 14.24773410052932
$ws.Range("L18").Value = 9.847697191994635
$ws.Range("M18").Value = 15.04182826674034
$ws.Range("O18").Value = 21.77766310466072
$ws.Range("B19").Value = 17.10211577299436
$ws.Range("C19").Value = 12.87233625940059
$ws.Range("D19").Value = 6.351762335751773
$ws.Range("E19").Value = 12.81856641907726
$ws.Range("G19").Value = 29.19898130270174
$ws.Range("H19").Value = 14.25195983319773
$ws.Range("L19").Value = 9.847775104954701
$ws.Range("M19").Value = 15.02953869532553
$ws.Range("O19").Value = 21.78350666775094
$ws.Range("B20").Value = 17.31536124019162
$ws.Range("C20").Value = 12.94239143809273
$ws.Range("D20").Value = 6.401474058566856
$ws.Range("E20").Value = 12.79859328667491
$ws.Range("G20").Value = 29.19375185123089
$ws.Range("H20").Value = 14.23314979502572
$ws.Range("L20").Value = 9.84751420806297
$ws.Range("M20").Value = 15.08487107910224
$ws.Range("O20").Value = 21.75764522815741
$ws.Range("B21").Value = 18.01388749138099
$ws.Range("C21").Value = 13.17491760032354
$ws.Range("D21").Value = 6.567586577949736
$ws.Range("E21").Value = 12.73358661833635
$ws.Range("G21").Value = 29.19168174961864
$ws.Range("H21").Value = 14.17356545435047
$ws.Range("L21").Value = 9.848204220758069
$ws.Range("M21").Value = 15.27114381847679
$ws.Range("O21").Value = 21.67839662525615
$ws.Range("B22").Value = 18.45647885962209
$ws.Range("C22").Value = 13.32458625731772
$ws.Range("D22").Value = 6.67526607469253
$ws.Range("E22").Value = 12.69266007333134
$ws.Range("G22").Value = 29.20192113444274
$ws.Range("H22").Value = 14.1373248090459
$ws.Range("L22").Value = 9.849818285244996
$ws.Range("M22").Value = 15.39300337705959
$ws.Range("O22").Value = 21.63232910609277
$ws.Range("B23").Value = 18.22156806616509
$ws.Range("C23").Value = 13.24492688139345
$ws.Range("D23").Value = 6.617889220201483
$ws.Range("E23").Value = 12.71436109022967
$ws.Range("G23").Value = 29.19539337839462
$ws.Range("H23").Value = 14.15641951797161
$ws.Range("L23").Value = 9.848850706761226
$ws.Range("M23").Value = 15.32796337225613
$ws.Range("O23").Value = 21.65639038399596
$ws.Range("B24").Value = 17.30371681933265
$ws.Range("C24").Value = 12.93855460225034
$ws.Range("D24").Value = 6.398747022271531
$ws.Range("E24").Value = 12.79968219800924
$ws.Range("G24").Value = 29.19398083554005
$ws.Range("H24").Value = 14.23416917123198
$ws.Range("L24").Value = 9.847522630775966
$ws.Range("M24").Value = 15.08183073419214
$ws.Range("O24").Value = 21.75903674264692
$ws.Range("B25").Value = 16.25765262718558
$ws.Range("C25").Value = 12.5993689236163
$ws.Range("D25").Value = 6.159924109541
$ws.Range("E25").Value = 12.89844988127349
$ws.Range("G25").Value = 29.24193345719794
$ws.Range("H25").Value = 14.32958659670582
$ws.Range("L25").Value = 9.851111782661505
$ws.Range("M25").Value = 14.81783268275709
$ws.Range("O25").Value = 21.89414309349839
